# Generate Report for Handoff
# Adds two new "Ready for handoff" file rows
#   8c5ee07d-11e9-4c55-96e3-c6771cfba2fd
#   961ad719-0757-46cc-818b-dfb6d6a18a06
# to the Overview / zh-cn / de-de worksheets, mirroring the existing rows.

$wb = $excel.ActiveWorkbook

$mdBase        = "https://github.com/OpenLocalizationTest/oltest/blob/bb9adedcefd87dae2a126fe63044560c813cb2c5/e2e"
$zhcnXlfBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a43aed920732cb51052db32aa5284c614a705d1f/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht"
$dedeXlfBase   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d8330e654a1fe382e5f80e2cb3870d1e92491102/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht"

$file1 = "8c5ee07d-11e9-4c55-96e3-c6771cfba2fd"
$hash1 = "48c559d5d82755f72a3d1f7acea9b42fa61662e7"
$file2 = "961ad719-0757-46cc-818b-dfb6d6a18a06"
$hash2 = "cf2c46e03f1b0aa4950cabf41c806d4d76dfb22b"

$status       = "Ready for handoff"
$handoffDate  = "2016-03-19 07:39:31"
$zhcnDatetime = "2016-03-19 07:39:23"
$dedeDatetime = $handoffDate
$epoch        = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "$file1.md"
$wsOverview.Range("B4").Value = $status
$wsOverview.Range("C4").Value = $status
$wsOverview.Range("D4").Value = $handoffDate
$wsOverview.Range("A4").Style = "HyperLink"
$wsOverview.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), "$mdBase/$file1.md", "", "", "$file1.md") | Out-Null

$wsOverview.Range("A5").Value = "$file2.md"
$wsOverview.Range("B5").Value = $status
$wsOverview.Range("C5").Value = $status
$wsOverview.Range("D5").Value = $handoffDate
$wsOverview.Range("A5").Style = "HyperLink"
$wsOverview.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsOverview.Hyperlinks.Add($wsOverview.Range("A5"), "$mdBase/$file2.md", "", "", "$file2.md") | Out-Null

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhXlf1 = "$file1.$hash1.zh-cn.xlf"
$zhXlf2 = "$file2.$hash2.zh-cn.xlf"

$wsZhCn.Range("A4").Value = "$file1.md"
$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = $status
$wsZhCn.Range("D4").Value = $zhXlf1
$wsZhCn.Range("E4").Value = $zhcnDatetime
$wsZhCn.Range("H4").Value = $epoch
$wsZhCn.Range("I4").Value = "'"
$wsZhCn.Range("J4").Value = "Include"
$wsZhCn.Range("A4").Style = "HyperLink"
$wsZhCn.Range("D4").Style = "HyperLink"
$wsZhCn.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "$mdBase/$file1.md", "", "", "$file1.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D4"), "$zhcnXlfBase/$zhXlf1", "", "", $zhXlf1) | Out-Null

$wsZhCn.Range("A5").Value = "$file2.md"
$wsZhCn.Range("B5").Value = ".md"
$wsZhCn.Range("C5").Value = $status
$wsZhCn.Range("D5").Value = $zhXlf2
$wsZhCn.Range("E5").Value = $zhcnDatetime
$wsZhCn.Range("H5").Value = $epoch
$wsZhCn.Range("I5").Value = "'"
$wsZhCn.Range("J5").Value = "Include"
$wsZhCn.Range("A5").Style = "HyperLink"
$wsZhCn.Range("D5").Style = "HyperLink"
$wsZhCn.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A5"), "$mdBase/$file2.md", "", "", "$file2.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("D5"), "$zhcnXlfBase/$zhXlf2", "", "", $zhXlf2) | Out-Null

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deXlf1 = "$file1.$hash1.de-de.xlf"
$deXlf2 = "$file2.$hash2.de-de.xlf"

$wsDeDe.Range("A4").Value = "$file1.md"
$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = $status
$wsDeDe.Range("D4").Value = $deXlf1
$wsDeDe.Range("E4").Value = $dedeDatetime
$wsDeDe.Range("H4").Value = $epoch
$wsDeDe.Range("I4").Value = "'"
$wsDeDe.Range("J4").Value = "Include"
$wsDeDe.Range("A4").Style = "HyperLink"
$wsDeDe.Range("D4").Style = "HyperLink"
$wsDeDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "$mdBase/$file1.md", "", "", "$file1.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D4"), "$dedeXlfBase/$deXlf1", "", "", $deXlf1) | Out-Null

$wsDeDe.Range("A5").Value = "$file2.md"
$wsDeDe.Range("B5").Value = ".md"
$wsDeDe.Range("C5").Value = $status
$wsDeDe.Range("D5").Value = $deXlf2
$wsDeDe.Range("E5").Value = $dedeDatetime
$wsDeDe.Range("H5").Value = $epoch
$wsDeDe.Range("I5").Value = "'"
$wsDeDe.Range("J5").Value = "Include"
$wsDeDe.Range("A5").Style = "HyperLink"
$wsDeDe.Range("D5").Style = "HyperLink"
$wsDeDe.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A5"), "$mdBase/$file2.md", "", "", "$file2.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("D5"), "$dedeXlfBase/$deXlf2", "", "", $deXlf2) | Out-Null
